# Student email is now pulled from the API instead of being maintained in
# this sheet, so the "Email" column (J) - including its mailto: hyperlinks
# and the dedicated "Hyperlink" cell style that only that column used - is
# removed. "Quê quán" (old column K) shifts left to become the new J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column J the way a user would before deleting it (this is what
# drives the sheet's recorded selection/active cell to J1:J1048576).
$ws.Columns("J").Select() | Out-Null

# Drop the mailto: hyperlinks that lived on J2:J4 ...
$ws.Hyperlinks.Delete() | Out-Null

# ... and the now-unused built-in "Hyperlink" cell style those cells used.
$wb.Styles.Item("Hyperlink").Delete() | Out-Null

# Finally remove the whole column; remaining columns (K -> J, etc.) shift left.
$ws.Columns("J").Delete() | Out-Null
